$wb = $excel.ActiveWorkbook

# --- PIR sheet: append new rows ---
$ws = $wb.Worksheets.Item("PIR")
$ws.Cells.Item(205,1).Value = "'2026-02-06"
$ws.Cells.Item(205,2).Value = "09:57:04"
$ws.Cells.Item(205,3).Value = "09:00"
$ws.Cells.Item(205,4).Value = "Bathroom"
$ws.Cells.Item(205,5).Value = "No Motion"
$ws.Cells.Item(205,6).Value = "Inactive"
$ws.Cells.Item(206,1).Value = "'2026-02-06"
$ws.Cells.Item(206,2).Value = "09:57:06"
$ws.Cells.Item(206,3).Value = "09:00"
$ws.Cells.Item(206,4).Value = "Bathroom"
$ws.Cells.Item(206,5).Value = "No Motion"
$ws.Cells.Item(206,6).Value = "Inactive"
$ws.Cells.Item(207,1).Value = "'2026-02-06"
$ws.Cells.Item(207,2).Value = "09:57:10"
$ws.Cells.Item(207,3).Value = "09:00"
$ws.Cells.Item(207,4).Value = "Bathroom"
$ws.Cells.Item(207,5).Value = "No Motion"
$ws.Cells.Item(207,6).Value = "Inactive"
$ws.Cells.Item(208,1).Value = "'2026-02-06"
$ws.Cells.Item(208,2).Value = "09:57:15"
$ws.Cells.Item(208,3).Value = "09:00"
$ws.Cells.Item(208,4).Value = "Bathroom"
$ws.Cells.Item(208,5).Value = "No Motion"
$ws.Cells.Item(208,6).Value = "Inactive"
$ws.Cells.Item(209,1).Value = "'2026-02-06"
$ws.Cells.Item(209,2).Value = "09:57:20"
$ws.Cells.Item(209,3).Value = "09:00"
$ws.Cells.Item(209,4).Value = "Bathroom"
$ws.Cells.Item(209,5).Value = "No Motion"
$ws.Cells.Item(209,6).Value = "Inactive"
$ws.Cells.Item(210,1).Value = "'2026-02-06"
$ws.Cells.Item(210,2).Value = "09:57:25"
$ws.Cells.Item(210,3).Value = "09:00"
$ws.Cells.Item(210,4).Value = "Bathroom"
$ws.Cells.Item(210,5).Value = "No Motion"
$ws.Cells.Item(210,6).Value = "Inactive"
$ws.Cells.Item(211,1).Value = "'2026-02-06"
$ws.Cells.Item(211,2).Value = "09:57:30"
$ws.Cells.Item(211,3).Value = "09:00"
$ws.Cells.Item(211,4).Value = "Bathroom"
$ws.Cells.Item(211,5).Value = "No Motion"
$ws.Cells.Item(211,6).Value = "Inactive"
$ws.Cells.Item(212,1).Value = "'2026-02-06"
$ws.Cells.Item(212,2).Value = "09:57:35"
$ws.Cells.Item(212,3).Value = "09:00"
$ws.Cells.Item(212,4).Value = "Bathroom"
$ws.Cells.Item(212,5).Value = "No Motion"
$ws.Cells.Item(212,6).Value = "Inactive"
$ws.Cells.Item(213,1).Value = "'2026-02-06"
$ws.Cells.Item(213,2).Value = "09:57:40"
$ws.Cells.Item(213,3).Value = "09:00"
$ws.Cells.Item(213,4).Value = "Bathroom"
$ws.Cells.Item(213,5).Value = "No Motion"
$ws.Cells.Item(213,6).Value = "Inactive"
$ws.Cells.Item(214,1).Value = "'2026-02-06"
$ws.Cells.Item(214,2).Value = "09:57:45"
$ws.Cells.Item(214,3).Value = "09:00"
$ws.Cells.Item(214,4).Value = "Bathroom"
$ws.Cells.Item(214,5).Value = "No Motion"
$ws.Cells.Item(214,6).Value = "Inactive"
$ws.Cells.Item(215,1).Value = "'2026-02-06"
$ws.Cells.Item(215,2).Value = "09:57:50"
$ws.Cells.Item(215,3).Value = "09:00"
$ws.Cells.Item(215,4).Value = "Bathroom"
$ws.Cells.Item(215,5).Value = "No Motion"
$ws.Cells.Item(215,6).Value = "Inactive"
$ws.Cells.Item(216,1).Value = "'2026-02-06"
$ws.Cells.Item(216,2).Value = "09:57:55"
$ws.Cells.Item(216,3).Value = "09:00"
$ws.Cells.Item(216,4).Value = "Bathroom"
$ws.Cells.Item(216,5).Value = "No Motion"
$ws.Cells.Item(216,6).Value = "Inactive"
$ws.Cells.Item(217,1).Value = "'2026-02-06"
$ws.Cells.Item(217,2).Value = "09:58:00"
$ws.Cells.Item(217,3).Value = "09:00"
$ws.Cells.Item(217,4).Value = "Bathroom"
$ws.Cells.Item(217,5).Value = "No Motion"
$ws.Cells.Item(217,6).Value = "Inactive"

# --- Humidity sheet: append new rows ---
$ws = $wb.Worksheets.Item("Humidity")
$ws.Cells.Item(122,1).Value = "'2026-02-06"
$ws.Cells.Item(122,2).Value = "09:57:04"
$ws.Cells.Item(122,3).Value = "09:00"
$ws.Cells.Item(122,4).Value = "Bathroom"
$ws.Cells.Item(122,5).Value = "'68.7%"
$ws.Cells.Item(122,6).Value = "Active"
$ws.Cells.Item(123,1).Value = "'2026-02-06"
$ws.Cells.Item(123,2).Value = "09:57:05"
$ws.Cells.Item(123,3).Value = "09:00"
$ws.Cells.Item(123,4).Value = "Bathroom"
$ws.Cells.Item(123,5).Value = "'69.1%"
$ws.Cells.Item(123,6).Value = "Active"
$ws.Cells.Item(124,1).Value = "'2026-02-06"
$ws.Cells.Item(124,2).Value = "09:57:24"
$ws.Cells.Item(124,3).Value = "09:00"
$ws.Cells.Item(124,4).Value = "Bathroom"
$ws.Cells.Item(124,5).Value = "'70.0%"
$ws.Cells.Item(124,6).Value = "Active"
$ws.Cells.Item(125,1).Value = "'2026-02-06"
$ws.Cells.Item(125,2).Value = "09:57:28"
$ws.Cells.Item(125,3).Value = "09:00"
$ws.Cells.Item(125,4).Value = "Bathroom"
$ws.Cells.Item(125,5).Value = "'70.0%"
$ws.Cells.Item(125,6).Value = "Active"
$ws.Cells.Item(126,1).Value = "'2026-02-06"
$ws.Cells.Item(126,2).Value = "09:57:33"
$ws.Cells.Item(126,3).Value = "09:00"
$ws.Cells.Item(126,4).Value = "Bathroom"
$ws.Cells.Item(126,5).Value = "'70.1%"
$ws.Cells.Item(126,6).Value = "Active"
$ws.Cells.Item(127,1).Value = "'2026-02-06"
$ws.Cells.Item(127,2).Value = "09:57:38"
$ws.Cells.Item(127,3).Value = "09:00"
$ws.Cells.Item(127,4).Value = "Bathroom"
$ws.Cells.Item(127,5).Value = "'70.2%"
$ws.Cells.Item(127,6).Value = "Active"
$ws.Cells.Item(128,1).Value = "'2026-02-06"
$ws.Cells.Item(128,2).Value = "09:57:43"
$ws.Cells.Item(128,3).Value = "09:00"
$ws.Cells.Item(128,4).Value = "Bathroom"
$ws.Cells.Item(128,5).Value = "'69.3%"
$ws.Cells.Item(128,6).Value = "Active"
$ws.Cells.Item(129,1).Value = "'2026-02-06"
$ws.Cells.Item(129,2).Value = "09:57:48"
$ws.Cells.Item(129,3).Value = "09:00"
$ws.Cells.Item(129,4).Value = "Bathroom"
$ws.Cells.Item(129,5).Value = "'70.2%"
$ws.Cells.Item(129,6).Value = "Active"
$ws.Cells.Item(130,1).Value = "'2026-02-06"
$ws.Cells.Item(130,2).Value = "09:57:53"
$ws.Cells.Item(130,3).Value = "09:00"
$ws.Cells.Item(130,4).Value = "Bathroom"
$ws.Cells.Item(130,5).Value = "'69.2%"
$ws.Cells.Item(130,6).Value = "Active"
$ws.Cells.Item(131,1).Value = "'2026-02-06"
$ws.Cells.Item(131,2).Value = "09:57:58"
$ws.Cells.Item(131,3).Value = "09:00"
$ws.Cells.Item(131,4).Value = "Bathroom"
$ws.Cells.Item(131,5).Value = "'70.2%"
$ws.Cells.Item(131,6).Value = "Active"
$ws.Cells.Item(132,1).Value = "'2026-02-06"
$ws.Cells.Item(132,2).Value = "09:58:03"
$ws.Cells.Item(132,3).Value = "09:00"
$ws.Cells.Item(132,4).Value = "Bathroom"
$ws.Cells.Item(132,5).Value = "'69.2%"
$ws.Cells.Item(132,6).Value = "Active"

# --- Temperature sheet: append new rows ---
$ws = $wb.Worksheets.Item("Temperature")
$ws.Cells.Item(122,1).Value = "'2026-02-06"
$ws.Cells.Item(122,2).Value = "09:57:04"
$ws.Cells.Item(122,3).Value = "09:00"
$ws.Cells.Item(122,4).Value = "Bathroom"
$ws.Cells.Item(122,5).Value = "27.9C"
$ws.Cells.Item(122,6).Value = "Active"
$ws.Cells.Item(123,1).Value = "'2026-02-06"
$ws.Cells.Item(123,2).Value = "09:57:05"
$ws.Cells.Item(123,3).Value = "09:00"
$ws.Cells.Item(123,4).Value = "Bathroom"
$ws.Cells.Item(123,5).Value = "27.8C"
$ws.Cells.Item(123,6).Value = "Active"
$ws.Cells.Item(124,1).Value = "'2026-02-06"
$ws.Cells.Item(124,2).Value = "09:57:24"
$ws.Cells.Item(124,3).Value = "09:00"
$ws.Cells.Item(124,4).Value = "Bathroom"
$ws.Cells.Item(124,5).Value = "27.9C"
$ws.Cells.Item(124,6).Value = "Active"
$ws.Cells.Item(125,1).Value = "'2026-02-06"
$ws.Cells.Item(125,2).Value = "09:57:28"
$ws.Cells.Item(125,3).Value = "09:00"
$ws.Cells.Item(125,4).Value = "Bathroom"
$ws.Cells.Item(125,5).Value = "27.9C"
$ws.Cells.Item(125,6).Value = "Active"
$ws.Cells.Item(126,1).Value = "'2026-02-06"
$ws.Cells.Item(126,2).Value = "09:57:33"
$ws.Cells.Item(126,3).Value = "09:00"
$ws.Cells.Item(126,4).Value = "Bathroom"
$ws.Cells.Item(126,5).Value = "27.9C"
$ws.Cells.Item(126,6).Value = "Active"
$ws.Cells.Item(127,1).Value = "'2026-02-06"
$ws.Cells.Item(127,2).Value = "09:57:38"
$ws.Cells.Item(127,3).Value = "09:00"
$ws.Cells.Item(127,4).Value = "Bathroom"
$ws.Cells.Item(127,5).Value = "27.9C"
$ws.Cells.Item(127,6).Value = "Active"
$ws.Cells.Item(128,1).Value = "'2026-02-06"
$ws.Cells.Item(128,2).Value = "09:57:44"
$ws.Cells.Item(128,3).Value = "09:00"
$ws.Cells.Item(128,4).Value = "Bathroom"
$ws.Cells.Item(128,5).Value = "27.9C"
$ws.Cells.Item(128,6).Value = "Active"
$ws.Cells.Item(129,1).Value = "'2026-02-06"
$ws.Cells.Item(129,2).Value = "09:57:49"
$ws.Cells.Item(129,3).Value = "09:00"
$ws.Cells.Item(129,4).Value = "Bathroom"
$ws.Cells.Item(129,5).Value = "27.9C"
$ws.Cells.Item(129,6).Value = "Active"
$ws.Cells.Item(130,1).Value = "'2026-02-06"
$ws.Cells.Item(130,2).Value = "09:57:54"
$ws.Cells.Item(130,3).Value = "09:00"
$ws.Cells.Item(130,4).Value = "Bathroom"
$ws.Cells.Item(130,5).Value = "27.9C"
$ws.Cells.Item(130,6).Value = "Active"
$ws.Cells.Item(131,1).Value = "'2026-02-06"
$ws.Cells.Item(131,2).Value = "09:57:59"
$ws.Cells.Item(131,3).Value = "09:00"
$ws.Cells.Item(131,4).Value = "Bathroom"
$ws.Cells.Item(131,5).Value = "27.9C"
$ws.Cells.Item(131,6).Value = "Active"
